$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '51.221.46'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.59%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.915.68'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.59%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.38%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '369.78'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +6.86%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '103.99'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -2.12%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.540'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -1.22%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.20%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.586'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -2.29%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.67'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -1.18%  '
$ws.Range("E11").Value = '  +1.38%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0835'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.73%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.35'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -1.81%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.373.58'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.17%  '
$ws.Range("E15").Value = '  -2.36%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.921.74'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.96%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.934'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -1.69%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '51.166.28'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.28%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.24'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -3.67%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.20'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -1.40%  '
$ws.Range("E21").Value = '  -2.07%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0944'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.69%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.41'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.21%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '259.61'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.34%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.68'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.90%  '
$ws.Range("E26").Value = '  +4.10%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.176'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +4.04%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.04%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '25.74'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.21%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.07'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -4.55%  '
$ws.Range("E31").Value = '  +0.03%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.22'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +4.08%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '9.89'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -1.84%  '
$ws.Range("E34").Value = '  +0.38%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '34.63'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.84%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '50.83'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.81%  '
$ws.Range("E37").Value = '  +0.05%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0424'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +1.45%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.01'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -1.78%  '
$ws.Range("E40").Value = '  +2.72%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '17.10'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -1.49%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.84'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -3.76%  '
$ws.Range("E43").Value = '  -1.17%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.18'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.89%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '118.50'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.21%  '
$ws.Range("E46").Value = '  -2.68%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.018.37'
$ws.Range("D47").ClearFormats()
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.31'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +2.81%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.17'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -2.75%  '
$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.213.26'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.52%  '
$ws.Range("B51").Value = 'TheGraph'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.241'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +2.35%  '
